# fix : player image with not exist -> make null image
# Replace the teamid values (column O, rows 2-24) with freshly generated
# identifiers, since the images referenced by the previous ids no longer
# exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTeamIds = @(
    "0874d67c-598f-44e3-920b-4561908dfe3b",
    "9e34f94d-c31f-4c23-99d2-a2822dd88d51",
    "003c4550-974d-4aa8-9276-e3ef21048e79",
    "fb8d99d6-9f17-469f-9dd5-003d5c7968d9",
    "9e9e4955-a3b5-41f2-a00c-bd1f5310c968",
    "1921564e-062f-4597-95dc-95074a120bb7",
    "a0768df2-6a70-45cb-9b17-170d678f0408",
    "bf18fea4-73ce-4d76-86f9-bbf70bcbff6a",
    "a75e188a-12c5-4285-bcec-2d97ca3c7153",
    "2111d74d-0d3e-4e4f-9f1b-b6e8c0cfe0f4",
    "50be74bd-92cf-49cc-97dc-5599f0670e91",
    "d55d68ed-1702-4f23-b775-3777dd67494c",
    "5c8943e1-fcd8-4718-b614-972b56a08435",
    "1c9e5a86-8cbe-46d6-989c-5e819a6aaf27",
    "5c8fd640-4c38-43fd-be97-b526717d2df2",
    "92c24ef1-b25c-4ff5-bb92-9058df7f7b6e",
    "2412507d-742b-4053-9f4e-e96905c184cd",
    "4c6a902c-c4c5-4ae2-88a1-05cd3ec101b8",
    "19995c9e-d80e-471a-9813-a42cd6f71119",
    "af854dff-40e7-4c50-a0f1-6dd6958c6784",
    "331efece-b4e3-4565-86c3-931d3f793263",
    "166314d0-c6c7-44e5-b2c2-2c69a7f25cd8",
    "32128616-e39a-4800-9f53-f18630bc75d0"
)

$startRow = 2
for ($i = 0; $i -lt $newTeamIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 15).Value = $newTeamIds[$i]
}
